# A new transaction log entry was recorded on the "2024" sheet's September
# column (R = September_Details, S = September_Date). This inserts a brand
# new row above the current most-recent entry (row 35), which pushes every
# row below it (35-110) down by one (to 36-111) and bumps the sheet's used
# range from A1:Y110 to A1:Y111. The new row gets the latest log entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new blank row before row 35 - shifts rows 35:110 down to 36:111.
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row with the latest September log entry.
$ws.Range("R35").Value = "bal axisbank"
$ws.Range("S35").Value = "2024-09-09 11:38:16"
